{"js": "// Replace the paragraph that describes the extra Employee fields / Customer\n// table with the updated wording from the commit.\nconst body = context.document.body;\n\nconst results = body.search(\"We added two additional fields to the Employee table\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst target = results.items[0];\nconst paras = target.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst paragraph = paras.items[0];\n\nconst newText =\n  \"For the extra thing that we thought would be useful, we added a new table labeled Customer. \" +\n  \"This table would allow the employees to also rate the customer and input a rating and a review \" +\n  \"in order keep track of people who may be barred from calling a taxi \\u2013 much like how Uber has a \" +\n  \"passenger rating and it allows drivers to deny customers if they have a low rating and a history of bad behaviour.\";\n\nparagraph.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph describing the extra Employee fields / Customer table.\n$range = $d.Content\n[void]$range.Find.Execute(\"We added two additional fields to the Employee table\")\n\n$paragraph = $range.Paragraphs(1).Range\n# Exclude the trailing paragraph mark so the whole paragraph body (and only\n# the body) gets replaced, while formatting of the existing run carries over.\n[void]$paragraph.MoveEnd(1, -1)\n\n$newText = \"For the extra thing that we thought would be useful, we added a new table labeled Customer. \" + `\n  \"This table would allow the employees to also rate the customer and input a rating and a review \" + `\n  \"in order keep track of \" + `\n  \"people who may be barred from calling a taxi \" + [char]0x2013 + \" much like how Uber has a passenger rating and it allows drivers to deny customers if they have a low rating and a history of bad behaviour.\"\n\n$paragraph.Text = $newText\n"}
